# bug fixing for centerfit
# Add two new trial rows to the becExpType table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18: EvapDD1PumpTime
$ws.Range("A18").Value = "EvapDD1PumpTime"
$ws.Range("B18").Value = "An experiment at evaporation stage D."
$ws.Range("C18").Value = "TOP"
$ws.Range("D18").Value = "EvapDOdt1"
$ws.Range("E18").Value = 4
$ws.Range("F18").Value = "D1PumpTime"
$ws.Range("G18").Value = "DensityFit;AtomNumber;CenterFit"
$ws.Range("H18").Value = "LSR"
$ws.Range("I18").Value = "LF"
$ws.Range("J18").Value = "RandomPolarization"
$ws.Range("K18").Value = 8
$ws.Range("L18").Value = "BosonicGaussianFit1D"
$ws.Range("M18").Value = 30
$ws.Range("N18").Value = "LinearFit1D"

# Row 19: HfBecTof
$ws.Range("A19").Value = "HfBecTof"
$ws.Range("B19").Value = "An experiment at the high-field BEC stage."
$ws.Range("C19").Value = "TOP"
$ws.Range("D19").Value = "Bec"
$ws.Range("E19").Value = 4
$ws.Range("F19").Value = "TOF"
$ws.Range("G19").Value = "AtomNumber;CenterFit;Tof;DensityFit"
$ws.Range("H19").Value = "LSR"
$ws.Range("I19").Value = "HF"
$ws.Range("J19").Value = "StrongLight"
$ws.Range("K19").Value = 8
$ws.Range("L19").Value = "BosonicGaussianFit1D"
$ws.Range("M19").Value = 1
$ws.Range("N19").Value = "ParabolicFit1D"
